# Update TPM-derived values in Vegfa-Nrp1 worksheet (columns G:T, rows 2-17)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 16,14
$arr[0,0] = 2.162809
$arr[0,1] = 6.488427000000001
$arr[0,2] = 0.06755089002018773
$arr[0,3] = 0.06755089002018773
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 127.3992563333333
$arr[0,7] = 382.197769
$arr[0,8] = 0.4838549810199306
$arr[0,9] = 0.4838549810199307
$arr[0,10] = 275.5402581910403
$arr[0,11] = 2479.862323719363
$arr[0,12] = 0.03268483460859736
$arr[0,13] = 0.03268483460859736
$arr[1,0] = 2.162809
$arr[1,1] = 6.488427000000001
$arr[1,2] = 0.06755089002018773
$arr[1,3] = 0.06755089002018773
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 59.36586533333332
$arr[1,7] = 178.097596
$arr[1,8] = 0.2254681108101269
$arr[1,9] = 0.2254681108101269
$arr[1,10] = 128.3970278357213
$arr[1,11] = 1155.573250521492
$arr[1,12] = 0.01523057155639438
$arr[1,13] = 0.01523057155639438
$arr[2,0] = 2.162809
$arr[2,1] = 6.488427000000001
$arr[2,2] = 0.06755089002018773
$arr[2,3] = 0.06755089002018773
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 16.63275166666667
$arr[2,7] = 49.898255
$arr[2,8] = 0.06317022542837675
$arr[2,9] = 0.06317022542837675
$arr[2,10] = 35.97346499943167
$arr[2,11] = 323.761184994885
$arr[2,12] = 0.004267204950462744
$arr[2,13] = 0.004267204950462744
$arr[3,0] = 2.162809
$arr[3,1] = 6.488427000000001
$arr[3,2] = 0.06755089002018773
$arr[3,3] = 0.06755089002018773
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 59.90262233333334
$arr[3,7] = 179.707867
$arr[3,8] = 0.2275066827415657
$arr[3,9] = 0.2275066827415658
$arr[3,10] = 129.5579307061344
$arr[3,11] = 1166.021376355209
$arr[3,12] = 0.01536827890473325
$arr[3,13] = 0.01536827890473325
$arr[4,0] = 18.019504
$arr[4,1] = 54.058512
$arr[4,2] = 0.5628021396814664
$arr[4,3] = 0.5628021396814664
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 127.3992563333333
$arr[4,7] = 382.197769
$arr[4,8] = 0.4838549810199306
$arr[4,9] = 0.4838549810199307
$arr[4,10] = 2295.671409095525
$arr[4,11] = 20661.04268185973
$arr[4,12] = 0.2723146186135523
$arr[4,13] = 0.2723146186135523
$arr[5,0] = 18.019504
$arr[5,1] = 54.058512
$arr[5,2] = 0.5628021396814664
$arr[5,3] = 0.5628021396814664
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 59.36586533333332
$arr[5,7] = 178.097596
$arr[5,8] = 0.2254681108101269
$arr[5,9] = 0.2254681108101269
$arr[5,10] = 1069.743447837461
$arr[5,11] = 9627.691030537151
$arr[5,12] = 0.1268939351938774
$arr[5,13] = 0.1268939351938774
$arr[6,0] = 18.019504
$arr[6,1] = 54.058512
$arr[6,2] = 0.5628021396814664
$arr[6,3] = 0.5628021396814664
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 16.63275166666667
$arr[6,7] = 49.898255
$arr[6,8] = 0.06317022542837675
$arr[6,9] = 0.06317022542837675
$arr[6,10] = 299.7139351885067
$arr[6,11] = 2697.42541669656
$arr[6,12] = 0.03555233803525101
$arr[6,13] = 0.03555233803525101
$arr[7,0] = 18.019504
$arr[7,1] = 54.058512
$arr[7,2] = 0.5628021396814664
$arr[7,3] = 0.5628021396814664
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 59.90262233333334
$arr[7,7] = 179.707867
$arr[7,8] = 0.2275066827415657
$arr[7,9] = 0.2275066827415658
$arr[7,10] = 1079.41554274599
$arr[7,11] = 9714.739884713905
$arr[7,12] = 0.1280412478387858
$arr[7,13] = 0.1280412478387858
$arr[8,0] = 4.650307000000001
$arr[8,1] = 13.950921
$arr[8,2] = 0.1452427730405732
$arr[8,3] = 0.1452427730405732
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 127.3992563333333
$arr[8,7] = 382.197769
$arr[8,8] = 0.4838549810199306
$arr[8,9] = 0.4838549810199307
$arr[8,10] = 592.4456535216943
$arr[8,11] = 5332.01088169525
$arr[8,12] = 0.07027643919282865
$arr[8,13] = 0.07027643919282867
$arr[9,0] = 4.650307000000001
$arr[9,1] = 13.950921
$arr[9,2] = 0.1452427730405732
$arr[9,3] = 0.1452427730405732
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 59.36586533333332
$arr[9,7] = 178.097596
$arr[9,8] = 0.2254681108101269
$arr[9,9] = 0.2254681108101269
$arr[9,10] = 276.0694991206573
$arr[9,11] = 2484.625492085916
$arr[9,12] = 0.03274761364628207
$arr[9,13] = 0.03274761364628208
$arr[10,0] = 4.650307000000001
$arr[10,1] = 13.950921
$arr[10,2] = 0.1452427730405732
$arr[10,3] = 0.1452427730405732
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 16.63275166666667
$arr[10,7] = 49.898255
$arr[10,8] = 0.06317022542837675
$arr[10,9] = 0.06317022542837675
$arr[10,10] = 77.34740150476168
$arr[10,11] = 696.126613542855
$arr[10,12] = 0.009175018714815571
$arr[10,13] = 0.009175018714815571
$arr[11,0] = 4.650307000000001
$arr[11,1] = 13.950921
$arr[11,2] = 0.1452427730405732
$arr[11,3] = 0.1452427730405732
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 59.90262233333334
$arr[11,7] = 179.707867
$arr[11,8] = 0.2275066827415657
$arr[11,9] = 0.2275066827415658
$arr[11,10] = 278.5655839550564
$arr[11,11] = 2507.090255595508
$arr[11,12] = 0.03304370148664693
$arr[11,13] = 0.03304370148664694
$arr[12,0] = 7.184856000000001
$arr[12,1] = 21.554568
$arr[12,2] = 0.2244041972577726
$arr[12,3] = 0.2244041972577726
$arr[12,4] = 3
$arr[12,5] = 1
$arr[12,6] = 127.3992563333333
$arr[12,7] = 382.197769
$arr[12,8] = 0.4838549810199306
$arr[12,9] = 0.4838549810199307
$arr[12,10] = 915.3453112620881
$arr[12,11] = 8238.107801358794
$arr[12,12] = 0.1085790886049523
$arr[12,13] = 0.1085790886049524
$arr[13,0] = 7.184856000000001
$arr[13,1] = 21.554568
$arr[13,2] = 0.2244041972577726
$arr[13,3] = 0.2244041972577726
$arr[13,4] = 3
$arr[13,5] = 1
$arr[13,6] = 59.36586533333332
$arr[13,7] = 178.097596
$arr[13,8] = 0.2254681108101269
$arr[13,9] = 0.2254681108101269
$arr[13,10] = 426.535193735392
$arr[13,11] = 3838.816743618528
$arr[13,12] = 0.05059599041357304
$arr[13,13] = 0.05059599041357305
$arr[14,0] = 7.184856000000001
$arr[14,1] = 21.554568
$arr[14,2] = 0.2244041972577726
$arr[14,3] = 0.2244041972577726
$arr[14,4] = 3
$arr[14,5] = 1
$arr[14,6] = 16.63275166666667
$arr[14,7] = 49.898255
$arr[14,8] = 0.06317022542837675
$arr[14,9] = 0.06317022542837675
$arr[14,10] = 119.50392560876
$arr[14,11] = 1075.53533047884
$arr[14,12] = 0.01417566372784742
$arr[14,13] = 0.01417566372784742
$arr[15,0] = 7.184856000000001
$arr[15,1] = 21.554568
$arr[15,2] = 0.2244041972577726
$arr[15,3] = 0.2244041972577726
$arr[15,4] = 3
$arr[15,5] = 1
$arr[15,6] = 59.90262233333334
$arr[15,7] = 179.707867
$arr[15,8] = 0.2275066827415657
$arr[15,9] = 0.2275066827415658
$arr[15,10] = 430.3917154873841
$arr[15,11] = 3873.525439386457
$arr[15,12] = 0.05105345451139981
$arr[15,13] = 0.05105345451139982

$ws.Range("G2:T17").Value = $arr
